$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the MOLD/Jan value (B4) from 4200000 to 3700000
$ws.Range("B4").Value = 3700000

# Move the active selection to B9 (reflecting where the user ended up after editing)
$ws.Range("B9").Select()

$wb.Save()
